# Generate Report for handoff
# Fill in the "Latest Handoff Datetime" (column D) for the c6d14bf7 file row (row 6)
# on both the "zh-cn" and "de-de" status sheets, which previously held a stale /
# duplicated timestamp copied from another row.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-01-15 07:46:39"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-01-15 07:46:49"
